$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

# Enterprises density (per 1000 people) - row 11
Set-TextValue "B11" "10.26"
Set-TextValue "C11" "0.81"
Set-TextValue "D11" "11.07"

# Employment (% of total) - row 12
Set-TextValue "B12" "51.27"
Set-TextValue "C12" "32.75"
Set-TextValue "D12" "84.02"

# Enterprises (% of total) - row 14
Set-TextValue "B14" "92.58"
Set-TextValue "C14" "7.34"
Set-TextValue "D14" "99.91"
